$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 332
$ws.Range("I33").Value = 273.5
$ws.Range("K33").Value = 273.5
$ws.Range("M33").Value = -44.5
$ws.Range("H62").Value = 2122.5
$ws.Range("I62").Value = 2122.5
$ws.Range("K62").Value = 2122.5
$ws.Range("M62").Value = -1498.5
$ws.Range("H65").Value = 2122.5
$ws.Range("I65").Value = 2122.5
$ws.Range("K65").Value = 10612.5
$ws.Range("M65").Value = -7492.5
$ws.Range("H74").Value = 6573.2
$ws.Range("J74").Value = 7676
$ws.Range("L74").Value = 7676
$ws.Range("N74").Value = -9548
$ws.Range("H77").Value = 6573.2
$ws.Range("J77").Value = 7676
$ws.Range("L77").Value = 38380
$ws.Range("N77").Value = -47740
$ws.Range("H93").Value = 42266.668
$ws.Range("J93").Value = 42266.668
$ws.Range("L93").Value = 42266.668
$ws.Range("N93").Value = -47258.668
$ws.Range("H137").Value = 655129.7
$ws.Range("I137").Value = 1255645.4
$ws.Range("J137").Value = 3141.257
$ws.Range("K137").Value = 3766936.2
$ws.Range("L137").Value = 9423.771000000001
$ws.Range("M137").Value = -3764386.2
$ws.Range("N137").Value = -14523.771
$ws.Range("H141").Value = 33036.363
$ws.Range("I141").Value = 47400
$ws.Range("J141").Value = 4309.091
$ws.Range("K141").Value = 142200
$ws.Range("L141").Value = 12927.273
$ws.Range("M141").Value = -137020
$ws.Range("N141").Value = -23287.273
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3232.051
$ws.Range("I32").Value = 3083.738
$ws.Range("K32").Value = 3083.738
$ws.Range("M32").Value = -2796.738
$ws.Range("H61").Value = 3320.5
$ws.Range("I61").Value = 3320.5
$ws.Range("K61").Value = 3320.5
$ws.Range("M61").Value = -3108.5
$ws.Range("H76").Value = 37600
$ws.Range("J76").Value = 37600
$ws.Range("L76").Value = 37600
$ws.Range("N76").Value = -38276
$ws.Range("H79").Value = 37600
$ws.Range("J79").Value = 37600
$ws.Range("L79").Value = 37600
$ws.Range("N79").Value = -39940
$ws.Range("H103").Value = 34400
$ws.Range("J103").Value = 34400
$ws.Range("L103").Value = 34400
$ws.Range("N103").Value = -36744
$ws.Range("H136").Value = 3320.5
$ws.Range("I136").Value = 3320.5
$ws.Range("K136").Value = 9961.5
$ws.Range("M136").Value = -7411.5
$ws.Range("H137").Value = 38269.715
$ws.Range("J137").Value = 38269.715
$ws.Range("L137").Value = 38269.715
$ws.Range("N137").Value = -48469.715
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H43").Value = 74800
$ws.Range("J43").Value = 74800
$ws.Range("L43").Value = 74800
$ws.Range("N43").Value = -75162
$ws.Range("H99").Value = 2667.1428
$ws.Range("I99").Value = 1382.5
$ws.Range("J99").Value = 4380
$ws.Range("K99").Value = 1382.5
$ws.Range("L99").Value = 4380
$ws.Range("M99").Value = 115.5
$ws.Range("N99").Value = -7376
$ws.Range("H107").Value = 1334.95
$ws.Range("I107").Value = 1327.1
$ws.Range("J107").Value = 1342.8
$ws.Range("K107").Value = 1327.1
$ws.Range("L107").Value = 1342.8
$ws.Range("M107").Value = 592.9000000000001
$ws.Range("N107").Value = -5182.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 216652.22
$ws.Range("I31").Value = 712108.5600000001
$ws.Range("J31").Value = 2705.1592
$ws.Range("K31").Value = 712108.5600000001
$ws.Range("L31").Value = 2705.1592
$ws.Range("M31").Value = -711813.5600000001
$ws.Range("N31").Value = -3295.1592
$ws.Range("H34").Value = 216652.22
$ws.Range("I34").Value = 712108.5600000001
$ws.Range("J34").Value = 2705.1592
$ws.Range("K34").Value = 712108.5600000001
$ws.Range("L34").Value = 2705.1592
$ws.Range("M34").Value = -711906.5600000001
$ws.Range("N34").Value = -3109.1592
$ws.Range("H107").Value = 3030978.5
$ws.Range("J107").Value = 1742.375
$ws.Range("L107").Value = 1742.375
$ws.Range("N107").Value = -5582.375
$ws.Range("H132").Value = 2412.25
$ws.Range("I132").Value = 1287.7727
$ws.Range("J132").Value = 6535.3335
$ws.Range("K132").Value = 3863.3181
$ws.Range("L132").Value = 19606.0005
$ws.Range("M132").Value = -1333.3181
$ws.Range("N132").Value = -24666.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4598.0586
$ws.Range("I68").Value = 1222
$ws.Range("J68").Value = 7599
$ws.Range("K68").Value = 3666
$ws.Range("L68").Value = 22797
$ws.Range("M68").Value = -2855
$ws.Range("N68").Value = -24419
$ws.Range("H71").Value = 4598.0586
$ws.Range("I71").Value = 1222
$ws.Range("J71").Value = 7599
$ws.Range("K71").Value = 10998
$ws.Range("L71").Value = 68391
$ws.Range("M71").Value = -6942
$ws.Range("N71").Value = -76503
$ws.Range("H104").Value = 5780
$ws.Range("J104").Value = 5780
$ws.Range("L104").Value = 17340
$ws.Range("N104").Value = -22582
$ws.Range("H107").Value = 14308.92
$ws.Range("I107").Value = 391.84616
$ws.Range("J107").Value = 29385.75
$ws.Range("K107").Value = 1175.53848
$ws.Range("L107").Value = 88157.25
$ws.Range("M107").Value = 744.4615200000001
$ws.Range("N107").Value = -91997.25
$ws.Range("H113").Value = 2084031.4
$ws.Range("I113").Value = 627.38635
$ws.Range("K113").Value = 1882.15905
$ws.Range("M113").Value = 287.84095
$ws.Range("H131").Value = 825.9299999999999
$ws.Range("I131").Value = 529
$ws.Range("J131").Value = 844.883
$ws.Range("K131").Value = 1587
$ws.Range("L131").Value = 2534.649
$ws.Range("M131").Value = 3453
$ws.Range("N131").Value = -12614.649
$ws.Range("H140").Value = 1887.5
$ws.Range("I140").Value = 1000
$ws.Range("J140").Value = 2014.2858
$ws.Range("K140").Value = 3000
$ws.Range("L140").Value = 6042.857400000001
$ws.Range("M140").Value = 2180
$ws.Range("N140").Value = -16402.8574
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1733.8889
$ws.Range("I113").Value = 1965.3334
$ws.Range("J113").Value = 1271
$ws.Range("K113").Value = 1965.3334
$ws.Range("L113").Value = 1271
$ws.Range("M113").Value = 204.6666
$ws.Range("N113").Value = -5611
$ws.Range("H141").Value = 32546.555
$ws.Range("J141").Value = 32546.555
$ws.Range("L141").Value = 32546.555
$ws.Range("N141").Value = -42906.555
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 60394.117
$ws.Range("I22").Value = 84033.336
$ws.Range("J22").Value = 3660
$ws.Range("K22").Value = 84033.336
$ws.Range("L22").Value = 3660
$ws.Range("M22").Value = -83738.336
$ws.Range("N22").Value = -4250
$ws.Range("H27").Value = 60394.117
$ws.Range("I27").Value = 84033.336
$ws.Range("J27").Value = 3660
$ws.Range("K27").Value = 84033.336
$ws.Range("L27").Value = 3660
$ws.Range("M27").Value = -83926.336
$ws.Range("N27").Value = -3874
$ws.Range("H40").Value = 5834.5835
$ws.Range("I40").Value = 4501.5
$ws.Range("J40").Value = 12500
$ws.Range("K40").Value = 4501.5
$ws.Range("L40").Value = 12500
$ws.Range("M40").Value = -4365.5
$ws.Range("N40").Value = -12772
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 33353334
$ws.Range("I49").Value = 100000000
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 100000000
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -99999770
$ws.Range("N49").Value = -30460
$ws.Range("H54").Value = 20504.6
$ws.Range("J54").Value = 20504.6
$ws.Range("L54").Value = 20504.6
$ws.Range("N54").Value = -21544.6
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = $null
$ws.Range("H113").Value = 438.5
$ws.Range("I113").Value = 451.33334
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 1354.00002
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 815.9999800000001
$ws.Range("N113").Value = -5540
$ws.Range("H136").Value = 2705.2046
$ws.Range("I136").Value = 1213
$ws.Range("J136").Value = 4495.85
$ws.Range("K136").Value = 3639
$ws.Range("L136").Value = 13487.55
$ws.Range("M136").Value = -1089
$ws.Range("N136").Value = -18587.55
